# #5: property aircraft done
# Fix property_category values: "建物" (building) sheet should say "building",
# and "汽車" (car) sheet should say "car" instead of the incorrect "land".

$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"
$wsBuilding.Range("I4").Value = "building"

$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
$wsCar.Range("H3").Value = "car"
